$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5 so the current totals row (row 5) moves down to row 6
$ws.Rows.Item(5).Insert()

# New row 5: duplicate "proprietaire" (mediexpets) info into new contrat "Ahmed Test"
$ws.Range("A5").Value = "Ahmed Test"
$ws.Range("B5").Value = "BG1949"
$ws.Range("D5").Value = "bmce"
$ws.Range("E5").Value = "bmce"
$ws.Range("F5").Value = "Supervision"
$ws.Range("G5").Value = "040/SUP SUD"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 15000
$ws.Range("J5").Value = 2250
$ws.Range("K5").Value = 12750

# Recalculate totals row (now row 6) to include the new row
$ws.Range("I6").Value = 155000
$ws.Range("J6").Value = 3050
$ws.Range("K6").Value = 151950
